# Report 13 02 2025 - apply price-history update to Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Correct the timestamp for the existing 2025-02-12 batch
#    (rows 218-261): 45700.8478157995 -> 45700.8478157986
# ------------------------------------------------------------------
for ($r = 218; $r -le 261; $r++) {
    $ws.Cells.Item($r, 1).Value = 45700.8478157986
}

# ------------------------------------------------------------------
# 2) Append the new 2025-02-13 batch of 22 rows (262-283)
# ------------------------------------------------------------------
$newRows = @(
  @{ Row=262; A=45701.8929438426; B='1 килограм сребърно кюлче Valcambi'; C=2034.37; D=2594.76; E='https://tavex.bg/srebro/1-kg-valcambi-sreburno-kiulche/' },
  @{ Row=263; A=45701.8929438426; B='1/10 унция сребърна монета Британия'; C=9.08; D=12.12; E='https://tavex.bg/srebro/1-10-unciya-srebarna-moneta-britaniya/' },
  @{ Row=264; A=45701.8929438426; B='1 унция сребърна монета Виенска Филхармония'; C=62.37; D=79.25; E='https://tavex.bg/srebro/1-unciq-srebyrna-avstriiska-filharmonia/' },
  @{ Row=265; A=45701.8929438426; B='1 унция сребърна монета канадски кленов лист'; C=62.37; D=79.98; E='https://tavex.bg/srebro/1-unciya-sreburen-kanadski-klenov-list/' },
  @{ Row=266; A=45701.8929438426; B='1 унция сребърна монета Австралийско кенгуру'; C=62.37; D=79.25; E='https://tavex.bg/srebro/1-oz-sreburna-moneta-avstraliysko-kenguru/' },
  @{ Row=267; A=45701.8929438426; B='1 унция сребърна монета Американски орел'; C=66.61; D=116.33; E='https://tavex.bg/srebro/1-oz-sreburen-orel/' },
  @{ Row=268; A=45701.8929438426; B='1 унция сребърен австралийски лунар Змия 2025'; C=78.72; D=174.5; E='https://tavex.bg/srebro/1-unciya-srebaren-avstraliyski-lunar-godina-na-zmiyata-2025/' },
  @{ Row=269; A=45701.8929438426; B='30 грама сребърна монета Китайска панда 2025'; C=64.24; D=154.28; E='https://tavex.bg/srebro/30-grama-srebarna-moneta-kitaiska-panda-2025/' },
  @{ Row=270; A=45701.8929438426; B='1 унция сребърен австралийски лунар Дракон 2024'; C=78.72; D=174.5; E='https://tavex.bg/srebro/1-unciya-srebyren-avstraliiski-lunar-drakon-2024/' },
  @{ Row=271; A=45701.8929438426; B='1 унция сребърен австралийски лунар Заек 2023'; C=78.72; D=203.58; E='https://tavex.bg/srebro/1-unciya-srebyren-avstraliiski-lunar-zaek-2023/' },
  @{ Row=272; A=45701.8929438426; B='30 грама сребърна монета Китайска панда 2024'; C=78.24; E='https://tavex.bg/srebro/30-grama-srebarna-moneta-kitayska-panda-2024/' },
  @{ Row=273; A=45701.8929438426; B='30 грама сребърна монета Китайска панда 2023'; C=78.24; E='https://tavex.bg/srebro/30-grama-srebarna-moneta-kitayska-panda-2023/' },
  @{ Row=274; A=45701.8929438426; B='1 унция сребърна монета Британия'; C=72.66; E='https://tavex.bg/srebro/1-unciya-srebarna-moneta-britaniya-2/' },
  @{ Row=275; A=45701.8929438426; B='1 унция Сребърна монета Кругерранд, Южна Африка'; C=72.66; E='https://tavex.bg/srebro/1-unciya-srebarna-moneta-krugerrand-yuzhna-afrika/' },
  @{ Row=276; A=45701.8929438426; B='25 бр. 1 унция сребърна монета Британия'; C=78.5; E='https://tavex.bg/srebro/25-broya-1-unciya-srebarna-moneta-britania-tubus/' },
  @{ Row=277; A=45701.8929438426; B='25 бр. 1 унция сребърна монета Кругерранд, Южна Африка'; C=78.5; E='https://tavex.bg/srebro/25-broya-1-unciya-srebarna-moneta-krugerrand-yujna-afrika/' },
  @{ Row=278; A=45701.8929438426; B='500 бр. 1 унция Сребърна монета Британия'; C=78.0; E='https://tavex.bg/srebro/500-broya-1-unciya-srebarna-moneta-britaniya-masterbox-kutiya/' },
  @{ Row=279; A=45701.8929438426; B='500 бр. 1 унция сребърна монета Кругерранд, Южна Африка'; C=78.0; E='https://tavex.bg/srebro/500-broya-1-unciya-srebarna-moneta-krugerrand-yujna-afrika/' },
  @{ Row=280; A=45701.8929438426; B='1 унция сребърна австралийска коала'; C=78.61; E='https://tavex.bg/srebro/1-unciya-srebyrna-avstraliiska-koala/' },
  @{ Row=281; A=45701.8929438426; B='1 унция сребърна монета австралийски лунар година на Тигъра 2022'; C=78.61; E='https://tavex.bg/srebro/1-unciya-srebyrna-moneta-avstraliiski-lunar-tigyr-2023/' },
  @{ Row=282; A=45701.8929438426; B='30 грама сребърна монета Китайска панда 2022'; C=78.24; E='https://tavex.bg/srebro/30-grama-srebyrna-kitayska-panda-2022/' },
  @{ Row=283; A=45701.8929438426; B='1 унция  Кукабура 2022 година'; C=78.61; E='https://tavex.bg/srebro/1-oz-australian-kookaburra-2022-silver-coin/' },
)

# Give the new rows the same number format as the existing timestamp
# column (yyyy-mm-dd hh:mm:ss) before writing the values, so column A
# keeps style index reused rather than minting a fresh one.
$ws.Range("A218").Copy() | Out-Null
$ws.Range("A262:A283").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

foreach ($row in $newRows) {
    $r = $row.Row
    $ws.Cells.Item($r, 1).Value = $row.A
    $ws.Cells.Item($r, 2).Value = $row.B
    $ws.Cells.Item($r, 3).Value = $row.C
    if ($row.ContainsKey('D')) {
        $ws.Cells.Item($r, 4).Value = $row.D
    }
    $ws.Cells.Item($r, 5).Value = $row.E
}

# ------------------------------------------------------------------
# 3) Drop the stale trailing placeholder rows that used to pad the
#    sheet out to the bottom of the worksheet.
# ------------------------------------------------------------------
$ws.Range("A1048568:E1048576").EntireRow.Delete() | Out-Null

# ------------------------------------------------------------------
# 4) Restore the view state (scroll position + active selection)
# ------------------------------------------------------------------
$ws.Range("C276").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 250
$excel.ActiveWindow.ScrollColumn = 1
